$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 1253.3392
$ws.Range("J112").Value = 1253.3392
$ws.Range("L112").Value = 3760.0176
$ws.Range("N112").Value = -5976.017599999999
# Row 135
$ws.Range("H135").Value = 83335740
$ws.Range("I135").Value = 33335832
$ws.Range("J135").Value = 333335300
$ws.Range("K135").Value = 300022488
$ws.Range("L135").Value = 3000017700
$ws.Range("M135").Value = -300019953
$ws.Range("N135").Value = -3000022770
# Row 141
$ws.Range("H141").Value = 4499.048
$ws.Range("I141").Value = 4335.625
$ws.Range("J141").Value = 5022
$ws.Range("K141").Value = 13006.875
$ws.Range("L141").Value = 15066
$ws.Range("M141").Value = -7826.875
$ws.Range("N141").Value = -25426

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 6382.9033
$ws.Range("I74").Value = 3807.4285
$ws.Range("J74").Value = 11791.4
$ws.Range("K74").Value = 3807.4285
$ws.Range("L74").Value = 11791.4
$ws.Range("M74").Value = -2933.4285
$ws.Range("N74").Value = -13539.4
# Row 77
$ws.Range("H77").Value = 6382.9033
$ws.Range("I77").Value = 3807.4285
$ws.Range("J77").Value = 11791.4
$ws.Range("K77").Value = 19037.1425
$ws.Range("L77").Value = 58957
$ws.Range("M77").Value = -14669.1425
$ws.Range("N77").Value = -67693
# Row 92
$ws.Range("H92").Value = 32250
$ws.Range("J92").Value = 32250
$ws.Range("L92").Value = 32250
$ws.Range("N92").Value = -37242

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1037.1818
$ws.Range("I20").Value = 1019.3125
$ws.Range("J20").Value = 1084.8334
$ws.Range("K20").Value = 1019.3125
$ws.Range("L20").Value = 1084.8334
$ws.Range("M20").Value = -772.3125
$ws.Range("N20").Value = -1578.8334
# Row 94
$ws.Range("H94").Value = 1726.7
$ws.Range("I94").Value = 1726.7
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1726.7
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1778
$ws.Range("I16").Value = 1599.5
$ws.Range("J16").Value = 1956.5
$ws.Range("K16").Value = 1599.5
$ws.Range("L16").Value = 1956.5
$ws.Range("M16").Value = -1312.5
$ws.Range("N16").Value = -2530.5
# Row 31
$ws.Range("H31").Value = 1879.3684
$ws.Range("I31").Value = 1413.3188
$ws.Range("K31").Value = 1413.3188
$ws.Range("M31").Value = -1118.3188
# Row 34
$ws.Range("H34").Value = 1879.3684
$ws.Range("I34").Value = 1413.3188
$ws.Range("K34").Value = 1413.3188
$ws.Range("M34").Value = -1211.3188
# Row 56
$ws.Range("H56").Value = 12000
$ws.Range("I56").Value = 12000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 12000
$ws.Range("M56").Value = -11155
$ws.Range("N56").ClearContents()
# Row 99
$ws.Range("H99").Value = 2852.75
$ws.Range("I99").Value = 2974.5715
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2974.5715
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -1476.5715
$ws.Range("N99").Value = -4996
# Row 107
$ws.Range("H107").Value = 973.3
$ws.Range("I107").Value = 958.55554
$ws.Range("J107").Value = 1106
$ws.Range("K107").Value = 958.55554
$ws.Range("L107").Value = 1106
$ws.Range("M107").Value = 961.44446
$ws.Range("N107").Value = -4946
# Row 113
$ws.Range("H113").Value = 1778
$ws.Range("I113").Value = 1599.5
$ws.Range("J113").Value = 1956.5
$ws.Range("K113").Value = 1599.5
$ws.Range("L113").Value = 1956.5
$ws.Range("M113").Value = 570.5
$ws.Range("N113").Value = -6296.5
# Row 126
$ws.Range("H126").Value = 2852.75
$ws.Range("I126").Value = 2974.5715
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 8923.7145
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -6453.7145
$ws.Range("N126").Value = -10940
# Row 132
$ws.Range("H132").Value = 2430.9565
$ws.Range("I132").Value = 2233.9524
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 6701.8572
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -4171.8572
$ws.Range("N132").Value = -18558.5
# Row 134
$ws.Range("H134").Value = 23323.295
$ws.Range("I134").Value = 41158.652
$ws.Range("J134").Value = 4774.52
$ws.Range("K134").Value = 123475.956
$ws.Range("L134").Value = 14323.56
$ws.Range("M134").Value = -120940.956
$ws.Range("N134").Value = -19393.56

$ws = $wb.Worksheets.Item("CUL")
# Row 64
$ws.Range("H64").Value = 2814.75
$ws.Range("I64").Value = 1322.8
$ws.Range("J64").Value = 3139.087
$ws.Range("K64").Value = 3968.4
$ws.Range("L64").Value = 9417.261
$ws.Range("M64").Value = -3698.4
$ws.Range("N64").Value = -9957.261
# Row 67
$ws.Range("H67").Value = 2814.75
$ws.Range("I67").Value = 1322.8
$ws.Range("J67").Value = 3139.087
$ws.Range("K67").Value = 3968.4
$ws.Range("L67").Value = 9417.261
$ws.Range("M67").Value = -3032.4
$ws.Range("N67").Value = -11289.261
# Row 70
$ws.Range("H70").Value = 2688.6
$ws.Range("I70").Value = 1148.2222
$ws.Range("J70").Value = 4999.1665
$ws.Range("K70").Value = 3444.6666
$ws.Range("L70").Value = 14997.4995
$ws.Range("M70").Value = -3129.6666
$ws.Range("N70").Value = -15627.4995
# Row 73
$ws.Range("H73").Value = 2688.6
$ws.Range("I73").Value = 1148.2222
$ws.Range("J73").Value = 4999.1665
$ws.Range("K73").Value = 3444.6666
$ws.Range("L73").Value = 14997.4995
$ws.Range("M73").Value = -2352.6666
$ws.Range("N73").Value = -17181.4995
# Row 75
$ws.Range("H75").Value = 333.33334
$ws.Range("I75").Value = 333.33334
$ws.Range("K75").Value = 1000.00002
$ws.Range("M75").Value = -2.000020000000063
# Row 78
$ws.Range("H78").Value = 333.33334
$ws.Range("I78").Value = 333.33334
$ws.Range("K78").Value = 3000.00006
$ws.Range("M78").Value = 1991.99994
# Row 87
$ws.Range("H87").Value = 7856.3
$ws.Range("I87").Value = 3465.6667
$ws.Range("J87").Value = 9738
$ws.Range("K87").Value = 10397.0001
$ws.Range("L87").Value = 29214
$ws.Range("M87").Value = -9149.000100000001
$ws.Range("N87").Value = -31710
# Row 90
$ws.Range("H90").Value = 7856.3
$ws.Range("I90").Value = 3465.6667
$ws.Range("J90").Value = 9738
$ws.Range("K90").Value = 31191.0003
$ws.Range("L90").Value = 87642
$ws.Range("M90").Value = -24951.0003
$ws.Range("N90").Value = -100122
# Row 122
$ws.Range("H122").Value = 705.2
$ws.Range("I122").Value = 393.89474
$ws.Range("J122").Value = 1074.875
$ws.Range("K122").Value = 3545.05266
$ws.Range("L122").Value = 9673.875
$ws.Range("M122").Value = -1095.05266
$ws.Range("N122").Value = -14573.875

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6599.05
$ws.Range("I70").Value = 6424.643
$ws.Range("J70").Value = 7006
$ws.Range("K70").Value = 6424.643
$ws.Range("L70").Value = 7006
$ws.Range("M70").Value = -6154.643
$ws.Range("N70").Value = -7546
# Row 73
$ws.Range("H73").Value = 6599.05
$ws.Range("I73").Value = 6424.643
$ws.Range("J73").Value = 7006
$ws.Range("K73").Value = 6424.643
$ws.Range("L73").Value = 7006
$ws.Range("M73").Value = -5488.643
$ws.Range("N73").Value = -8878

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3134.875
$ws.Range("I7").Value = 2276
$ws.Range("J7").Value = 3993.75
$ws.Range("K7").Value = 2276
$ws.Range("L7").Value = 3993.75
$ws.Range("M7").Value = -2164
$ws.Range("N7").Value = -4217.75
# Row 126
$ws.Range("H126").Value = 3134.875
$ws.Range("I126").Value = 2276
$ws.Range("J126").Value = 3993.75
$ws.Range("K126").Value = 6828
$ws.Range("L126").Value = 11981.25
$ws.Range("M126").Value = -4358
$ws.Range("N126").Value = -16921.25

$ws = $wb.Worksheets.Item("WVR")
# Row 80
$ws.Range("H80").Value = 33575.25
$ws.Range("J80").Value = 34767
$ws.Range("L80").Value = 34767
$ws.Range("N80").Value = -36763
# Row 83
$ws.Range("H83").Value = 33575.25
$ws.Range("J83").Value = 34767
$ws.Range("L83").Value = 104301
$ws.Range("N83").Value = -114285
# Row 122
$ws.Range("H122").Value = 2309.3
$ws.Range("I122").Value = 2272.9473
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6818.841899999999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4368.841899999999
$ws.Range("N122").Value = -13900
# Row 126
$ws.Range("H126").Value = 1141
$ws.Range("I126").Value = 886.75
$ws.Range("J126").Value = 1649.5
$ws.Range("K126").Value = 2660.25
$ws.Range("L126").Value = 4948.5
$ws.Range("M126").Value = -190.25
$ws.Range("N126").Value = -9888.5
# Row 132
$ws.Range("H132").Value = 1590.7188
$ws.Range("I132").Value = 1329.7778
$ws.Range("K132").Value = 3989.3334
$ws.Range("M132").Value = -1459.3334
# Row 136
$ws.Range("H136").Value = 7012
$ws.Range("I136").Value = 3438
$ws.Range("J136").Value = 9871.200000000001
$ws.Range("K136").Value = 10314
$ws.Range("L136").Value = 29613.6
$ws.Range("M136").Value = -7764
$ws.Range("N136").Value = -34713.60000000001
